$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3109
$ws.Range("I74").Value = 3181.6667
$ws.Range("J74").Value = 3000
$ws.Range("K74").Value = 3181.6667
$ws.Range("L74").Value = 3000
$ws.Range("M74").Value = -2245.6667
$ws.Range("N74").Value = -4872
$ws.Range("H77").Value = 3109
$ws.Range("I77").Value = 3181.6667
$ws.Range("J77").Value = 3000
$ws.Range("K77").Value = 15908.3335
$ws.Range("L77").Value = 15000
$ws.Range("M77").Value = -11228.3335
$ws.Range("N77").Value = -24360
$ws.Range("H80").Value = 1440.5454
$ws.Range("I80").Value = 2481.3
$ws.Range("J80").Value = 573.25
$ws.Range("K80").Value = 7443.900000000001
$ws.Range("L80").Value = 1719.75
$ws.Range("M80").Value = -6445.900000000001
$ws.Range("N80").Value = -3715.75
$ws.Range("H83").Value = 1440.5454
$ws.Range("I83").Value = 2481.3
$ws.Range("J83").Value = 573.25
$ws.Range("K83").Value = 22331.7
$ws.Range("L83").Value = 5159.25
$ws.Range("M83").Value = -17339.7
$ws.Range("N83").Value = -15143.25
$ws.Range("H112").Value = 5184.3076
$ws.Range("J112").Value = 5733.8
$ws.Range("L112").Value = 17201.4
$ws.Range("N112").Value = -19417.4
$ws.Range("H137").Value = 1543.0227
$ws.Range("I137").Value = 1487.1818
$ws.Range("J137").Value = 1598.8636
$ws.Range("K137").Value = 4461.5454
$ws.Range("L137").Value = 4796.5908
$ws.Range("M137").Value = -1911.5454
$ws.Range("N137").Value = -9896.5908
$ws.Range("H138").Value = 2542.6667
$ws.Range("I138").Value = 1366.7222
$ws.Range("J138").Value = 2803.9875
$ws.Range("K138").Value = 4100.1666
$ws.Range("L138").Value = 8411.962500000001
$ws.Range("M138").Value = 1039.8334
$ws.Range("N138").Value = -18691.9625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5943.25
$ws.Range("I32").Value = 5067.287
$ws.Range("J32").Value = 19666.666
$ws.Range("K32").Value = 5067.287
$ws.Range("L32").Value = 19666.666
$ws.Range("M32").Value = -4780.287
$ws.Range("N32").Value = -20240.666
$ws.Range("H61").Value = 7248681.5
$ws.Range("I61").Value = 10102610
$ws.Range("K61").Value = 10102610
$ws.Range("M61").Value = -10102398
$ws.Range("H63").Value = 49431.25
$ws.Range("I63").Value = 452150
$ws.Range("J63").Value = 4684.722
$ws.Range("K63").Value = 452150
$ws.Range("L63").Value = 4684.722
$ws.Range("M63").Value = -451464
$ws.Range("N63").Value = -6056.722
$ws.Range("H66").Value = 49431.25
$ws.Range("I66").Value = 452150
$ws.Range("J66").Value = 4684.722
$ws.Range("K66").Value = 2260750
$ws.Range("L66").Value = 23423.61
$ws.Range("M66").Value = -2257318
$ws.Range("N66").Value = -30287.61
$ws.Range("H136").Value = 7248681.5
$ws.Range("I136").Value = 10102610
$ws.Range("K136").Value = 30307830
$ws.Range("M136").Value = -30305280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1312.3334
$ws.Range("I99").Value = 1171.4286
$ws.Range("J99").Value = 1805.5
$ws.Range("K99").Value = 1171.4286
$ws.Range("L99").Value = 1805.5
$ws.Range("M99").Value = 326.5714
$ws.Range("N99").Value = -4801.5
$ws.Range("H105").Value = 2257.1428
$ws.Range("I105").Value = 2257.1428
$ws.Range("K105").Value = 2257.1428
$ws.Range("M105").Value = -510.1428000000001
$ws.Range("H106").Value = 76557
$ws.Range("J106").Value = 76557
$ws.Range("L106").Value = 76557
$ws.Range("N106").Value = -79081
$ws.Range("H114").Value = 66500
$ws.Range("J114").Value = 66500
$ws.Range("L114").Value = 66500
$ws.Range("N114").Value = -75178

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5127.9536
$ws.Range("I31").Value = 1829.9546
$ws.Range("J31").Value = 6815.3022
$ws.Range("K31").Value = 1829.9546
$ws.Range("L31").Value = 6815.3022
$ws.Range("M31").Value = -1534.9546
$ws.Range("N31").Value = -7405.3022
$ws.Range("H34").Value = 5127.9536
$ws.Range("I34").Value = 1829.9546
$ws.Range("J34").Value = 6815.3022
$ws.Range("K34").Value = 1829.9546
$ws.Range("L34").Value = 6815.3022
$ws.Range("M34").Value = -1627.9546
$ws.Range("N34").Value = -7219.3022
$ws.Range("H109").Value = 30285
$ws.Range("J109").Value = 30285
$ws.Range("L109").Value = 30285
$ws.Range("N109").Value = -32365
$ws.Range("H132").Value = 21741780
$ws.Range("I132").Value = 35716564
$ws.Range("J132").Value = 3223.4443
$ws.Range("K132").Value = 107149692
$ws.Range("L132").Value = 9670.332900000001
$ws.Range("M132").Value = -107147162
$ws.Range("N132").Value = -14730.3329

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 2116.2632
$ws.Range("J39").Value = 2116.2632
$ws.Range("L39").Value = 6348.7896
$ws.Range("N39").Value = -6936.7896
$ws.Range("H93").Value = 4977.353
$ws.Range("J93").Value = 4980.8125
$ws.Range("L93").Value = 14942.4375
$ws.Range("N93").Value = -18686.4375
$ws.Range("H107").Value = 2625.3333
$ws.Range("J107").Value = 3659.4
$ws.Range("L107").Value = 10978.2
$ws.Range("N107").Value = -14818.2
$ws.Range("H108").Value = 1551
$ws.Range("I108").Value = 1084.6666
$ws.Range("K108").Value = 3253.9998
$ws.Range("M108").Value = -373.9998000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 7562.5
$ws.Range("J15").Value = 7562.5
$ws.Range("L15").Value = 7562.5
$ws.Range("N15").Value = -8138.5
$ws.Range("H62").Value = 31076.111
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("H65").Value = 31076.111
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("H81").Value = 7562.5
$ws.Range("J81").Value = 7562.5
$ws.Range("L81").Value = 7562.5
$ws.Range("N81").Value = -9558.5
$ws.Range("H84").Value = 7562.5
$ws.Range("J84").Value = 7562.5
$ws.Range("L84").Value = 22687.5
$ws.Range("N84").Value = -32671.5
$ws.Range("H102").Value = 2162.0908
$ws.Range("I102").Value = 2128.3
$ws.Range("J102").Value = 2500
$ws.Range("K102").Value = 2128.3
$ws.Range("L102").Value = 2500
$ws.Range("M102").Value = -506.3000000000002
$ws.Range("N102").Value = -5744
$ws.Range("H113").Value = 93619.414
$ws.Range("I113").Value = 101857.55
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 101857.55
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = -99687.55
$ws.Range("N113").Value = -7340
$ws.Range("H132").Value = 22729558
$ws.Range("I132").Value = 30304814
$ws.Range("J132").Value = 3793.0908
$ws.Range("K132").Value = 90914442
$ws.Range("L132").Value = 11379.2724
$ws.Range("M132").Value = -90911912
$ws.Range("N132").Value = -16439.2724
$ws.Range("M62").ClearContents()
$ws.Range("M65").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2110.889
$ws.Range("I16").Value = 1533.1666
$ws.Range("J16").Value = 3266.3333
$ws.Range("K16").Value = 1533.1666
$ws.Range("L16").Value = 3266.3333
$ws.Range("M16").Value = -1363.1666
$ws.Range("N16").Value = -3606.3333
$ws.Range("H82").Value = 2075.1875
$ws.Range("I82").Value = 1716.6666
$ws.Range("J82").Value = 2290.3
$ws.Range("K82").Value = 1716.6666
$ws.Range("L82").Value = 2290.3
$ws.Range("M82").Value = -1355.6666
$ws.Range("N82").Value = -3012.3
$ws.Range("H85").Value = 2075.1875
$ws.Range("I85").Value = 1716.6666
$ws.Range("J85").Value = 2290.3
$ws.Range("K85").Value = 1716.6666
$ws.Range("L85").Value = 2290.3
$ws.Range("M85").Value = -468.6666
$ws.Range("N85").Value = -4786.3
$ws.Range("H132").Value = 3309.4075
$ws.Range("I132").Value = 2888.4443
$ws.Range("J132").Value = 3519.889
$ws.Range("K132").Value = 8665.332900000001
$ws.Range("L132").Value = 10559.667
$ws.Range("M132").Value = -6135.332900000001
$ws.Range("N132").Value = -15619.667
$ws.Range("H140").Value = 53687.727
$ws.Range("J140").Value = 53687.727
$ws.Range("L140").Value = 53687.727
$ws.Range("N140").Value = -64047.727

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value = 98325
$ws.Range("J86").Value = 98325
$ws.Range("L86").Value = 98325
$ws.Range("N86").Value = -100571
$ws.Range("H89").Value = 98325
$ws.Range("J89").Value = 98325
$ws.Range("L89").Value = 491625
$ws.Range("N89").Value = -502857
$ws.Range("H100").Value = 692.9048
$ws.Range("I100").Value = 599.5454999999999
$ws.Range("J100").Value = 795.6
$ws.Range("K100").Value = 1199.091
$ws.Range("L100").Value = 1591.2
$ws.Range("M100").Value = -658.0909999999999
$ws.Range("N100").Value = -2673.2
$ws.Range("H109").Value = 8000000
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("H122").Value = 2305.5334
$ws.Range("I122").Value = 2424.4167
$ws.Range("J122").Value = 1830
$ws.Range("K122").Value = 7273.250100000001
$ws.Range("L122").Value = 5490
$ws.Range("M122").Value = -4823.250100000001
$ws.Range("N122").Value = -10390
$ws.Range("N109").ClearContents()
